$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column B (EXECUTE) for rows 2-7 and 9 from "Yes" to "No"
$ws.Range("B2").Value = "No"
$ws.Range("B3").Value = "No"
$ws.Range("B4").Value = "No"
$ws.Range("B5").Value = "No"
$ws.Range("B6").Value = "No"
$ws.Range("B7").Value = "No"
$ws.Range("B9").Value = "No"

# Add new test row: A20 = "Test", B20 = "Yes"
$ws.Range("A20").Value = "Test"
$ws.Range("B20").Value = "Yes"

# Update selection to B11
$ws.Range("B11").Select()
